$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 31250300
$ws.Range("I2").Value = 31250300
$ws.Range("K2").Value = 31250300
$ws.Range("M2").Value = -31250187
$ws.Range("H17").Value = 3487875
$ws.Range("J17").Value = 3487875
$ws.Range("L17").Value = 10463625
$ws.Range("N17").Value = -10463961
$ws.Range("H106").Value = 63494012
$ws.Range("I106").Value = 37038748
$ws.Range("J106").Value = 83335460
$ws.Range("K106").Value = 37038748
$ws.Range("L106").Value = 83335460
$ws.Range("M106").Value = -37038117
$ws.Range("N106").Value = -83336722
$ws.Range("H107").Value = 6944998
$ws.Range("I107").Value = 10000310
$ws.Range("J107").Value = 1106.909
$ws.Range("K107").Value = 10000310
$ws.Range("L107").Value = 1106.909
$ws.Range("M107").Value = -9998390
$ws.Range("N107").Value = -4946.909
$ws.Range("H137").Value = 1531.16
$ws.Range("I137").Value = 1413.1875
$ws.Range("J137").Value = 1740.8889
$ws.Range("K137").Value = 4239.5625
$ws.Range("L137").Value = 5222.6667
$ws.Range("M137").Value = -1689.5625
$ws.Range("N137").Value = -10322.6667

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 14287089
$ws.Range("I74").Value = 1259.6522
$ws.Range("J74").Value = 41668260
$ws.Range("K74").Value = 1259.6522
$ws.Range("L74").Value = 41668260
$ws.Range("M74").Value = -385.6522
$ws.Range("N74").Value = -41670008
$ws.Range("H77").Value = 14287089
$ws.Range("I77").Value = 1259.6522
$ws.Range("J77").Value = 41668260
$ws.Range("K77").Value = 6298.261
$ws.Range("L77").Value = 208341300
$ws.Range("M77").Value = -1930.261
$ws.Range("N77").Value = -208350036
$ws.Range("H88").Value = 4333.3335
$ws.Range("J88").Value = 4333.3335
$ws.Range("L88").Value = 4333.3335
$ws.Range("N88").Value = -5145.3335
$ws.Range("H91").Value = 4333.3335
$ws.Range("J91").Value = 4333.3335
$ws.Range("L91").Value = 4333.3335
$ws.Range("N91").Value = -7141.3335
$ws.Range("H122").Value = 1166752.4
$ws.Range("I122").Value = 2332160
$ws.Range("J122").Value = 1344.8182
$ws.Range("K122").Value = 6996480
$ws.Range("L122").Value = 4034.4546
$ws.Range("M122").Value = -6994030
$ws.Range("N122").Value = -8934.454600000001
$ws.Range("H123").Value = 37708.5
$ws.Range("J123").Value = 37708.5
$ws.Range("L123").Value = 37708.5
$ws.Range("N123").Value = -47508.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 219.9
$ws.Range("I80").Value = 102
$ws.Range("J80").Value = 240.70589
$ws.Range("K80").Value = 102
$ws.Range("L80").Value = 240.70589
$ws.Range("M80").Value = 896
$ws.Range("N80").Value = -2236.70589
$ws.Range("H83").Value = 219.9
$ws.Range("I83").Value = 102
$ws.Range("J83").Value = 240.70589
$ws.Range("K83").Value = 510
$ws.Range("L83").Value = 1203.52945
$ws.Range("M83").Value = 4482
$ws.Range("N83").Value = -11187.52945
$ws.Range("H86").Value = 22224444
$ws.Range("I86").Value = 30304986
$ws.Range("J86").Value = 2951.75
$ws.Range("K86").Value = 30304986
$ws.Range("L86").Value = 2951.75
$ws.Range("M86").Value = -30303863
$ws.Range("N86").Value = -5197.75
$ws.Range("H89").Value = 22224444
$ws.Range("I89").Value = 30304986
$ws.Range("J89").Value = 2951.75
$ws.Range("K89").Value = 151524930
$ws.Range("L89").Value = 14758.75
$ws.Range("M89").Value = -151519314
$ws.Range("N89").Value = -25990.75
$ws.Range("H99").Value = 66667696
$ws.Range("I99").Value = 76924040
$ws.Range("K99").Value = 76924040
$ws.Range("M99").Value = -76922542

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15155337
$ws.Range("I31").Value = 1797
$ws.Range("J31").Value = 31255972
$ws.Range("K31").Value = 1797
$ws.Range("L31").Value = 31255972
$ws.Range("M31").Value = -1502
$ws.Range("N31").Value = -31256562
$ws.Range("H34").Value = 15155337
$ws.Range("I34").Value = 1797
$ws.Range("J34").Value = 31255972
$ws.Range("K34").Value = 1797
$ws.Range("L34").Value = 31255972
$ws.Range("M34").Value = -1595
$ws.Range("N34").Value = -31256376
$ws.Range("H122").Value = 2078.25
$ws.Range("I122").Value = 2362.4
$ws.Range("J122").Value = 1604.6666
$ws.Range("K122").Value = 7087.200000000001
$ws.Range("L122").Value = 4813.9998
$ws.Range("M122").Value = -4637.200000000001
$ws.Range("N122").Value = -9713.9998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7696.514
$ws.Range("I3").Value = 10667.538
$ws.Range("J3").Value = 5940.909
$ws.Range("K3").Value = 32002.614
$ws.Range("L3").Value = 17822.727
$ws.Range("M3").Value = -31890.614
$ws.Range("N3").Value = -18046.727
$ws.Range("H22").Value = 2244.6667
$ws.Range("J22").Value = 2671.7144
$ws.Range("L22").Value = 8015.1432
$ws.Range("N22").Value = -8353.143199999999
$ws.Range("J23").Value = 97.818184
$ws.Range("L23").Value = 293.454552
$ws.Range("N23").Value = -763.454552
$ws.Range("H27").Value = 2244.6667
$ws.Range("J27").Value = 2671.7144
$ws.Range("L27").Value = 8015.1432
$ws.Range("N27").Value = -8219.143199999999
$ws.Range("H32").Value = 3336.3635
$ws.Range("J32").Value = 3988.889
$ws.Range("L32").Value = 11966.667
$ws.Range("N32").Value = -12532.667
$ws.Range("H35").Value = 790
$ws.Range("J35").Value = 790
$ws.Range("L35").Value = 2370
$ws.Range("N35").Value = -2946
$ws.Range("H46").Value = 2576.923
$ws.Range("J46").Value = 2576.923
$ws.Range("L46").Value = 7730.768999999999
$ws.Range("N46").Value = -7912.768999999999
$ws.Range("H47").Value = 1225
$ws.Range("J47").Value = 2200
$ws.Range("L47").Value = 6600
$ws.Range("N47").Value = -7462
$ws.Range("H49").Value = 850
$ws.Range("J49").Value = 1375
$ws.Range("L49").Value = 4125
$ws.Range("N49").Value = -4437
$ws.Range("H54").Value = 3000
$ws.Range("J54").Value = 3000
$ws.Range("L54").Value = 9000
$ws.Range("N54").Value = -10118
$ws.Range("H58").Value = 2921.3684
$ws.Range("I58").Value = 1750
$ws.Range("J58").Value = 3059.1765
$ws.Range("K58").Value = 5250
$ws.Range("L58").Value = 9177.529500000001
$ws.Range("M58").Value = -5122
$ws.Range("N58").Value = -9433.529500000001
$ws.Range("H61").Value = 362.5
$ws.Range("J61").Value = 463.33334
$ws.Range("L61").Value = 1390.00002
$ws.Range("N61").Value = -1820.00002
$ws.Range("H80").Value = 550
$ws.Range("J80").Value = 600
$ws.Range("L80").Value = 1800
$ws.Range("N80").Value = -3672
$ws.Range("H83").Value = 550
$ws.Range("J83").Value = 600
$ws.Range("L83").Value = 5400
$ws.Range("N83").Value = -14760
$ws.Range("H94").Value = 2903.2
$ws.Range("I94").Value = 824
$ws.Range("J94").Value = 3223.077
$ws.Range("K94").Value = 2472
$ws.Range("L94").Value = 9669.231
$ws.Range("M94").Value = -1796
$ws.Range("N94").Value = -11021.231
$ws.Range("H105").Value = 6942.143
$ws.Range("J105").Value = 6942.143
$ws.Range("L105").Value = 20826.429
$ws.Range("N105").Value = -26068.429
$ws.Range("H106").Value = 3481.818
$ws.Range("J106").Value = 3481.818
$ws.Range("L106").Value = 10445.454
$ws.Range("N106").Value = -12337.454
$ws.Range("H113").Value = 2400497
$ws.Range("I113").Value = 2632076
$ws.Range("J113").Value = 1667164
$ws.Range("K113").Value = 7896228
$ws.Range("L113").Value = 5001492
$ws.Range("M113").Value = -7894058
$ws.Range("N113").Value = -5005832
$ws.Range("H115").Value = 1624.75
$ws.Range("I115").Value = 749.5
$ws.Range("K115").Value = 2248.5
$ws.Range("M115").Value = -1073.5
$ws.Range("H117").Value = 23814300
$ws.Range("I117").Value = 331.6
$ws.Range("J117").Value = 37044284
$ws.Range("K117").Value = 994.8000000000001
$ws.Range("L117").Value = 111132852
$ws.Range("M117").Value = 2447.2
$ws.Range("N117").Value = -111139736
$ws.Range("H131").Value = 3334218.2
$ws.Range("I131").Value = 11111607
$ws.Range("J131").Value = 1051.9048
$ws.Range("K131").Value = 33334821
$ws.Range("L131").Value = 3155.7144
$ws.Range("M131").Value = -33329781
$ws.Range("N131").Value = -13235.7144
$ws.Range("H138").Value = 6193.32
$ws.Range("I138").Value = 8655.799999999999
$ws.Range("K138").Value = 25967.4
$ws.Range("M138").Value = -20827.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5851.7393
$ws.Range("I70").Value = 5562.5
$ws.Range("K70").Value = 5562.5
$ws.Range("M70").Value = -5292.5
$ws.Range("H73").Value = 5851.7393
$ws.Range("I73").Value = 5562.5
$ws.Range("K73").Value = 5562.5
$ws.Range("M73").Value = -4626.5
$ws.Range("H132").Value = 11113948
$ws.Range("I132").Value = 18521110
$ws.Range("J132").Value = 3204.3333
$ws.Range("K132").Value = 55563330
$ws.Range("L132").Value = 9612.999899999999
$ws.Range("M132").Value = -55560800
$ws.Range("N132").Value = -14672.9999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 724.625
$ws.Range("I16").Value = 713.8570999999999
$ws.Range("K16").Value = 713.8570999999999
$ws.Range("M16").Value = -543.8570999999999
$ws.Range("H17").Value = 1700
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = -330
$ws.Range("N17").Value = -2340
$ws.Range("H132").Value = 13097929
$ws.Range("I132").Value = 16977724
$ws.Range("J132").Value = 3623.125
$ws.Range("K132").Value = 50933172
$ws.Range("L132").Value = 10869.375
$ws.Range("M132").Value = -50930642
$ws.Range("N132").Value = -15929.375

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2982.5
$ws.Range("J96").Value = 2935
$ws.Range("L96").Value = 2935
$ws.Range("N96").Value = -5681
$ws.Range("H123").Value = 34919
$ws.Range("J123").Value = 34919
$ws.Range("L123").Value = 34919
$ws.Range("N123").Value = -44719
$ws.Range("H132").Value = 1458.9459
$ws.Range("I132").Value = 1054.6
$ws.Range("J132").Value = 3191.8572
$ws.Range("K132").Value = 3163.8
$ws.Range("L132").Value = 9575.571599999999
$ws.Range("M132").Value = -633.7999999999997
$ws.Range("N132").Value = -14635.5716
